$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# StartStatus sheet: add Level column, shift CraftLevel, add AtkDelay/AtkSpeed
# columns, and expand the single data row into a 17-row progression table.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("StartStatus")

# Header row (row 2): ID, Level, InventorySlot, AddCoin, AtkPower, DefPower,
# Health, CraftLevel, AtkDelay, AtkSpeed
$ws.Cells.Item(2,1).Value = "ID"
$ws.Cells.Item(2,2).Value = "Level"
$ws.Cells.Item(2,3).Value = "InventorySlot"
$ws.Cells.Item(2,4).Value = "AddCoin"
$ws.Cells.Item(2,5).Value = "AtkPower"
$ws.Cells.Item(2,6).Value = "DefPower"
$ws.Cells.Item(2,7).Value = "Health"
$ws.Cells.Item(2,8).Value = "CraftLevel"
$ws.Cells.Item(2,9).Value = "AtkDelay"
$ws.Cells.Item(2,10).Value = "AtkSpeed"

# Data rows 3..19 (A..B, D..J are plain values; C is a formula column)
$startStatusRows = @(
    @(710000,0,5,50,2,100,1,0.2,0.5),
    @(710001,1,10,55,4,150,2,0.2,0.5),
    @(710002,2,15,60,6,200,3,0.2,0.5),
    @(710003,3,20,65,8,250,4,0.2,0.5),
    @(710004,4,25,70,10,300,5,0.2,0.5),
    @(710005,5,30,75,12,350,6,0.2,0.5),
    @(710006,6,35,80,14,400,7,0.2,0.5),
    @(710007,7,40,85,16,450,8,0.2,0.5),
    @(710008,8,45,90,18,500,9,0.2,0.5),
    @(710009,9,50,95,20,550,10,0.2,0.5),
    @(710010,10,55,100,22,600,11,0.2,0.5),
    @(710011,11,60,105,24,650,12,0.2,0.5),
    @(710012,12,65,110,26,700,13,0.2,0.5),
    @(710013,13,70,115,28,750,14,0.2,0.5),
    @(710014,14,75,120,30,800,15,0.2,0.5),
    @(710015,15,80,125,32,850,16,0.2,0.5),
    @(710016,16,85,130,34,900,17,0.2,0.5)
)

for ($i = 0; $i -lt $startStatusRows.Length; $i++) {
    $r = 3 + $i
    $row = $startStatusRows[$i]
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    if ($r -eq 3) {
        $ws.Cells.Item($r,3).Value = 4
    } else {
        $prev = $r - 1
        $ws.Cells.Item($r,3).Formula = "=IF(C$prev=20,20,SUM(C$prev,1))"
    }
    $ws.Cells.Item($r,4).Value = $row[2]
    $ws.Cells.Item($r,5).Value = $row[3]
    $ws.Cells.Item($r,6).Value = $row[4]
    $ws.Cells.Item($r,7).Value = $row[5]
    $ws.Cells.Item($r,8).Value = $row[6]
    $ws.Cells.Item($r,9).Value = $row[7]
    $ws.Cells.Item($r,10).Value = $row[8]
}

$ws.Range("K1:K1048576").Select()

# ---------------------------------------------------------------------------
# Upgrade sheet: add Level column, split each upgrade type into level 0/1
# rows with updated pricing.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Upgrade")

# Header row (row 2): ID, UpgradeType, Level, IconID, Price, IncreaseStat
$ws.Cells.Item(2,1).Value = "ID"
$ws.Cells.Item(2,2).Value = "UpgradeType"
$ws.Cells.Item(2,3).Value = "Level"
$ws.Cells.Item(2,4).Value = "IconID"
$ws.Cells.Item(2,5).Value = "Price"
$ws.Cells.Item(2,6).Value = "IncreaseStat"

$upgradeRows = @(
    @(401000,"Attack",0,940001,5000,5),
    @(401001,"Attack",1,940001,6000,5),
    @(402000,"Defence",0,940002,10000,1),
    @(402001,"Defence",1,940002,11000,1),
    @(403000,"Health",0,940003,2500,10),
    @(403001,"Health",1,940003,3000,10)
)

for ($i = 0; $i -lt $upgradeRows.Length; $i++) {
    $r = 3 + $i
    $row = $upgradeRows[$i]
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
}

$ws.Range("E9").Select()

# ---------------------------------------------------------------------------
# Monster sheet: selection only.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Monster")
$ws.Range("C33").Select()

# ---------------------------------------------------------------------------
# _Schedule sheet: remove V5, add V10, update selection (done last so this
# sheet stays the active tab, matching the source state).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("_Schedule")
$ws.Range("V5").ClearContents()
$ws.Range("V10").Value = 2
$ws.Range("W6").Select()
